# Apply the "add monte_carlo and update database" edit to the income
# statement workbook: refresh the published-date header for the latest
# period and fill in the previously-placeholder (all zero / "-") figures
# for the 1401/12-ended yearly column set (D:H across the statement rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 4 + $i   # D = 4 .. H = 8
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

# --- Row 9: "تاریخ انتشار" (publish date) header refresh for the last two periods ---
$ws.Cells.Item(9, 7).Value = "1402-03-07 (8)"
$ws.Cells.Item(9, 8).Value = "1402-03-07 (2)"

# --- Income statement figures (rows 11-27, columns D:H) ---

# فروش (Sales)
Set-RowValues 11 @(2591304, 3847004, 5921840, 7232077, 12355447)

# بهای تمام شده کالای فروش رفته (Cost of goods sold)
Set-RowValues 12 @(-1189529, -1834030, -2589100, -4184126, -6279464)

# سود (زیان) ناخالص (Gross profit)
Set-RowValues 13 @(1401775, 2012974, 3332740, 3047951, 6075983)

# هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
Set-RowValues 14 @(-269746, -176597, -193472, -437487, -654396)

# هزینه کاهش ارزش دریافتنی‌ها (هزینه استثنایی) - was placeholder "-" text, now 0
Set-RowValues 15 @(0, 0, 0, 0, 0)

# خالص سایر درامدها (هزینه ها) ی عملیاتی
Set-RowValues 16 @(121654, 2523, 66881, -25756, 75335)

# سود (زیان) عملیاتی (Operating profit)
Set-RowValues 17 @(1253683, 1838900, 3206149, 2584708, 5496922)

# هزینه های مالی (Finance costs)
Set-RowValues 18 @(-241803, -262692, -309456, -558612, -1011636)

# خالص سایر درامدها و هزینه های غیرعملیاتی
Set-RowValues 19 @(-12817, 33048, 79488, 141105, 91530)

# سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
Set-RowValues 20 @(999063, 1609256, 2976181, 2167201, 4576816)

# مالیات (Tax)
Set-RowValues 21 @(-234558, -357766, -487114, -362758, -750790)

# سود (زیان) خالص عملیات در حال تداوم
Set-RowValues 22 @(764505, 1251490, 2489067, 1804443, 3826026)

# سود (زیان) عملیات متوقف شده پس از اثر مالیاتی - was placeholder "-" text, now 0
Set-RowValues 23 @(0, 0, 0, 0, 0)

# سود (زیان) خالص (Net profit)
Set-RowValues 24 @(764505, 1251490, 2489067, 1804443, 3826026)

# سود هر سهم پس از کسر مالیات (EPS)
Set-RowValues 25 @(742, 1215, 1627, 743, 1054)

# سرمایه (Capital)
Set-RowValues 26 @(1030000, 1030000, 1530000, 2430000, 3630000)

# سود هر سهم بر اساس آخرین سرمایه
Set-RowValues 27 @(211, 345, 686, 497, 1054)

"done"
